$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 635.6
$ws.Range("I6").Value = 127
$ws.Range("K6").Value = 381
$ws.Range("M6").Value = -269
$ws.Range("H9").Value = 79.5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 79.5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 79.5
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = -417.5
$ws.Range("H17").Value = 814.04346
$ws.Range("J17").Value = 814.04346
$ws.Range("L17").Value = 2442.13038
$ws.Range("N17").Value = -2778.13038
$ws.Range("H97").Value = 1198.5
$ws.Range("J97").Value = 1198.5
$ws.Range("L97").Value = 3595.5
$ws.Range("N97").Value = -4587.5
$ws.Range("H100").Value = 1761.9
$ws.Range("I100").Value = 1624.3334
$ws.Range("K100").Value = 1624.3334
$ws.Range("M100").Value = -1083.3334
$ws.Range("H108").Value = 62000
$ws.Range("J108").Value = 62000
$ws.Range("L108").Value = 62000
$ws.Range("N108").Value = -69680
$ws.Range("H125").Value = 392
$ws.Range("I125").Value = 392
$ws.Range("K125").Value = 3528
$ws.Range("M125").Value = -1068
$ws.Range("H129").Value = 910.9524
$ws.Range("J129").Value = 881.5
$ws.Range("L129").Value = 2644.5
$ws.Range("N129").Value = -12644.5
$ws.Range("H132").Value = 1258.24
$ws.Range("I132").Value = 1185.6666
$ws.Range("K132").Value = 3556.9998
$ws.Range("M132").Value = -1026.9998
$ws.Range("H138").Value = 2830.4055
$ws.Range("I138").Value = 2454.276
$ws.Range("J138").Value = 4193.875
$ws.Range("K138").Value = 7362.828
$ws.Range("L138").Value = 12581.625
$ws.Range("M138").Value = -2222.828
$ws.Range("N138").Value = -22861.625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6301.8335
$ws.Range("I61").Value = 3270.3333
$ws.Range("K61").Value = 3270.3333
$ws.Range("M61").Value = -3058.3333
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988
$ws.Range("H110").Value = 1478.8462
$ws.Range("I110").Value = 792
$ws.Range("K110").Value = 792
$ws.Range("M110").Value = 1253
$ws.Range("H122").Value = 1700
$ws.Range("I122").Value = 1100
$ws.Range("K122").Value = 3300
$ws.Range("M122").Value = -850
$ws.Range("H123").Value = 53426
$ws.Range("J123").Value = 53426
$ws.Range("L123").Value = 53426
$ws.Range("N123").Value = -63226
$ws.Range("H132").Value = 2417.818
$ws.Range("I132").Value = 2055
$ws.Range("K132").Value = 6165
$ws.Range("M132").Value = -3635
$ws.Range("H136").Value = 6301.8335
$ws.Range("I136").Value = 3270.3333
$ws.Range("K136").Value = 9810.999899999999
$ws.Range("M136").Value = -7260.999899999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 299.67743
$ws.Range("I94").Value = 310
$ws.Range("J94").Value = 150
$ws.Range("K94").Value = 310
$ws.Range("L94").Value = 150
$ws.Range("M94").Value = 141
$ws.Range("N94").Value = -1052
$ws.Range("H134").Value = 7754.136
$ws.Range("I134").Value = 7754.136
$ws.Range("K134").Value = 23262.408
$ws.Range("M134").Value = -20727.408

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1838.6666
$ws.Range("I31").Value = 1146.1578
$ws.Range("K31").Value = 1146.1578
$ws.Range("M31").Value = -851.1578
$ws.Range("H34").Value = 1838.6666
$ws.Range("I34").Value = 1146.1578
$ws.Range("K34").Value = 1146.1578
$ws.Range("M34").Value = -944.1578
$ws.Range("H106").Value = 35335.5
$ws.Range("J106").Value = 35671
$ws.Range("L106").Value = 35671
$ws.Range("N106").Value = -38195
$ws.Range("H134").Value = 983.0769
$ws.Range("I134").Value = 965
$ws.Range("K134").Value = 2895
$ws.Range("M134").Value = -360
$ws.Range("H140").Value = 39499.5
$ws.Range("J140").Value = 39499.5
$ws.Range("L140").Value = 39499.5
$ws.Range("N140").Value = -49859.5
$ws.Range("H141").Value = 35764.3
$ws.Range("J141").Value = 61528.6
$ws.Range("L141").Value = 61528.6
$ws.Range("N141").Value = -71888.60000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2832.6667
$ws.Range("I69").Value = 2832.6667
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 8498.000100000001
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -7687.000100000001
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 2832.6667
$ws.Range("I72").Value = 2832.6667
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 25494.0003
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -21438.0003
$ws.Range("N72").Value = $null
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H107").Value = 767
$ws.Range("J107").Value = 798.8889
$ws.Range("L107").Value = 2396.6667
$ws.Range("N107").Value = -6236.6667
$ws.Range("H131").Value = 766.50507
$ws.Range("J131").Value = 782.6882000000001
$ws.Range("L131").Value = 2348.0646
$ws.Range("N131").Value = -12428.0646

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 59.5
$ws.Range("I2").Value = 27.2
$ws.Range("J2").Value = 74.181816
$ws.Range("K2").Value = 27.2
$ws.Range("L2").Value = 74.181816
$ws.Range("M2").Value = 85.8
$ws.Range("N2").Value = -300.181816
$ws.Range("H80").Value = 3322.2856
$ws.Range("I80").Value = 3269
$ws.Range("J80").Value = 3393.3333
$ws.Range("K80").Value = 3269
$ws.Range("L80").Value = 3393.3333
$ws.Range("M80").Value = -2271
$ws.Range("N80").Value = -5389.3333
$ws.Range("H83").Value = 3322.2856
$ws.Range("I83").Value = 3269
$ws.Range("J83").Value = 3393.3333
$ws.Range("K83").Value = 16345
$ws.Range("L83").Value = 16966.6665
$ws.Range("M83").Value = -11353
$ws.Range("N83").Value = -26950.6665
$ws.Range("H102").Value = 2337.9167
$ws.Range("I102").Value = 2295.3157
$ws.Range("J102").Value = 2499.8
$ws.Range("K102").Value = 2295.3157
$ws.Range("L102").Value = 2499.8
$ws.Range("M102").Value = -673.3157000000001
$ws.Range("N102").Value = -5743.8
$ws.Range("H113").Value = 1392.75
$ws.Range("I113").Value = 1098
$ws.Range("J113").Value = 1491
$ws.Range("K113").Value = 1098
$ws.Range("L113").Value = 1491
$ws.Range("M113").Value = 1072
$ws.Range("N113").Value = -5831
$ws.Range("H122").Value = 3002.1667
$ws.Range("I122").Value = 2006.5
$ws.Range("K122").Value = 6019.5
$ws.Range("M122").Value = -3569.5
$ws.Range("H126").Value = 1769781.8
$ws.Range("I126").Value = 1986574.1
$ws.Range("J126").Value = 252235
$ws.Range("K126").Value = 5959722.300000001
$ws.Range("L126").Value = 756705
$ws.Range("M126").Value = -5957252.300000001
$ws.Range("N126").Value = -761645

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2719.2
$ws.Range("I16").Value = 2642.923
$ws.Range("K16").Value = 2642.923
$ws.Range("M16").Value = -2472.923
$ws.Range("H68").Value = 2537.077
$ws.Range("I68").Value = 2368.9
$ws.Range("J68").Value = 3097.6667
$ws.Range("K68").Value = 2368.9
$ws.Range("L68").Value = 3097.6667
$ws.Range("M68").Value = -1619.9
$ws.Range("N68").Value = -4595.6667
$ws.Range("H71").Value = 2537.077
$ws.Range("I71").Value = 2368.9
$ws.Range("J71").Value = 3097.6667
$ws.Range("K71").Value = 11844.5
$ws.Range("L71").Value = 15488.3335
$ws.Range("M71").Value = -8100.5
$ws.Range("N71").Value = -22976.3335
$ws.Range("H88").Value = 49499.5
$ws.Range("J88").Value = 49499.5
$ws.Range("L88").Value = 49499.5
$ws.Range("N88").Value = -50355.5
$ws.Range("H91").Value = 49499.5
$ws.Range("J91").Value = 49499.5
$ws.Range("L91").Value = 49499.5
$ws.Range("N91").Value = -52463.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 45306.062
$ws.Range("J123").Value = 47499.785
$ws.Range("L123").Value = 47499.785
$ws.Range("N123").Value = -57299.785
$ws.Range("H132").Value = 1048.7142
$ws.Range("I132").Value = 802.5
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 2407.5
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = 122.5
$ws.Range("N132").Value = -17808.5
